$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TTNDY")

# Replace specific J-column (2018 period) values with "NA" for rows that
# no longer have reported data.
$ws.Range("J21").Value = "NA"
$ws.Range("J83").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"

# Update Capital Expenditures (row 91) figures across the period columns.
$ws.Range("D91").Value = -203600
$ws.Range("E91").Value = -186600
$ws.Range("F91").Value = -159800
$ws.Range("G91").Value = -144000
$ws.Range("H91").Value = -104600
$ws.Range("I91").Value = -102400
$ws.Range("J91").Value = -94600
